$wb = $excel.ActiveWorkbook

# Remove the empty "Feuil2" sheet, leaving "Feuil1" as the only sheet.
$excel.DisplayAlerts = $false
[void]$wb.Worksheets("Feuil2").Delete()

$ws = $wb.Worksheets("Feuil1")
$ws.Activate()

# Update the view (scrolled to row 6, normal zoom 100%).
$win = $excel.ActiveWindow
$win.Zoom = 100
$win.ScrollRow = 6
$win.ScrollColumn = 1

# Fix the wording in B17: drop the word "populaires".
$ws.Range("B17").Value = "Pas d'application officielle, mais des apps `ntierces comme Jerboa (Android) `net Memmy (iOS)"
